# Add 2022-Q3 data
# --------------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" right after "总计" and before
#    "2022-Q2" (copying "2022-Q2" so the new sheet inherits the same
#    header/column formatting), then populate it with the Q3 fund data.
# 2) Insert a new row at the top of the "总计" summary sheet for the
#    "2022-Q3" totals, pushing the older quarters down by one row.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet
# ---------------------------------------------------------------------
$qSheet = $wb.Worksheets.Item("2022-Q2")
$qSheet.Copy($qSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Clear any leftover rows from the copied sheet and make sure the text
# columns (B, C, D, E, F, G) keep their values verbatim (no silent
# numeric coercion / lost leading zeros / lost trailing zeros).
$newSheet.Range("A2:H9").ClearContents()
$newSheet.Range("B2:G9").NumberFormat = "@"

# Rows 8 and 9 are new (the copied sheet only had 7 rows) - give column A
# the same "ID column" style used by the rest of column A (copied from A7).
$newSheet.Range("A7").Copy()
$newSheet.Range("A8:A9").PasteSpecial(-4122)

$funds = @(
    @{ code = "012671"; name = "嘉实核心蓝筹混合A"; size = "9.20"; total = "93.58"; pct = "4.40"; mv = "0.4048"; rank = 9 },
    @{ code = "009126"; name = "嘉实基础产业优选股票A"; size = "1.74"; total = "93.85"; pct = "4.35"; mv = "0.0757"; rank = 8 },
    @{ code = "011924"; name = "嘉实港股互联网产业核心资产混合A"; size = "1.20"; total = "87.88"; pct = "4.22"; mv = "0.0506"; rank = 6 },
    @{ code = "012672"; name = "嘉实核心蓝筹混合C"; size = "0.42"; total = "93.58"; pct = "4.40"; mv = "0.0185"; rank = 9 },
    @{ code = "011925"; name = "嘉实港股互联网产业核心资产混合C"; size = "0.41"; total = "87.88"; pct = "4.22"; mv = "0.0173"; rank = 6 },
    @{ code = "013897"; name = "德邦港股通成长精选混合型证券投资基金A"; size = "0.41"; total = "79.99"; pct = "3.44"; mv = "0.0141"; rank = 7 },
    @{ code = "013898"; name = "德邦港股通成长精选混合型证券投资基金C"; size = "0.37"; total = "79.99"; pct = "3.44"; mv = "0.0127"; rank = 7 },
    @{ code = "009127"; name = "嘉实基础产业优选股票C"; size = "0.21"; total = "93.85"; pct = "4.35"; mv = "0.0091"; rank = 8 }
)

$row = 2
foreach ($f in $funds) {
    $idx = $row - 2
    $newSheet.Range("A$row").Value = $idx
    $newSheet.Range("B$row").Value = $f.code
    $newSheet.Range("C$row").Value = $f.name
    $newSheet.Range("D$row").Value = $f.size
    $newSheet.Range("E$row").Value = $f.total
    $newSheet.Range("F$row").Value = $f.pct
    $newSheet.Range("G$row").Value = $f.mv
    $newSheet.Range("H$row").Value = $f.rank
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Step 2: add the "2022-Q3" row to the "总计" (totals) summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $totalSheet.Range("A$dst").Value = $totalSheet.Range("A$src").Value()
    $totalSheet.Range("B$dst").Value = $totalSheet.Range("B$src").Value()
    $totalSheet.Range("C$dst").Value = $totalSheet.Range("C$src").Value()
    $totalSheet.Range("D$dst").Value = $totalSheet.Range("D$src").Value()
}

# New row 9 needs the same "ID column" style as the rest of column A.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 0.6

$newSheet.Activate()
